# Refresh the "cryptos" price/volume snapshot (Price = column D, Volume(1h) = column E).
# Note: several Price values (e.g. "0.4450", "313.28") look numeric, so a bare
# assignment would make Excel silently coerce them to numbers (dropping
# trailing zeros / losing the original text representation). Prefixing with
# a leading apostrophe forces them to stay literal text, matching the source
# workbook, and the cell Style is reset to "Normal" afterwards so it doesn't
# pick up the implicit "Text" number format that the apostrophe entry mode
# would otherwise leave behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.409.05"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "1.827.65"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'313.28"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4450"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "'0.3780"
$ws.Range("E8").Value = "  +3.01%  "
$ws.Range("D9").Value = "'0.07397"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").Value = "'0.8799"
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "1.833.77"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "'6.727"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").Value = "'5.434"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("D15").Value = "'92.92"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "'0.07058"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "27.415.07"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").Value = "'5.366"
$ws.Range("E22").Value = "  +4.44%  "
$ws.Range("D23").Value = "'10.96"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("D24").Value = "'1.948"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'150.98"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'2.287"
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("D27").Value = "'18.64"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").Value = "'5.358"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("D29").Value = "'117.19"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").Value = "'0.08907"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").Value = "'0.7922"
$ws.Range("E31").Value = "  +6.27%  "
$ws.Range("D32").Value = "'1.200"
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "'2.930"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.9997"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("D37").Value = "'0.01980"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "'7.306"
$ws.Range("E39").Value = "  +3.67%  "
$ws.Range("D40").Value = "'0.5325"
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").Value = "'2.874"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "'2.344"
$ws.Range("E42").Value = "  +19.01%  "
$ws.Range("D43").Value = "'0.1701"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "'8.679"
$ws.Range("E44").Value = "  +3.05%  "
$ws.Range("D45").Value = "'0.5065"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "'10.64"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").Value = "'105.55"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "'0.9994"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'0.06389"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").Value = "'66.39"
$ws.Range("E51").Value = "  +6.13%  "

# Reset number format/style back to Normal for cells forced to text via apostrophe
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
